{"js": "// Helper: wrap a <w:body> inner fragment into a full OOXML \"flat package\"\n// document suitable for Range.insertOoxml().\nfunction wrapBody(inner) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + inner + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items.forEach(p => p.load(\"text\"));\nawait context.sync();\n\n// Locate every paragraph we need to touch by its distinctive text.\nlet idxRunTest = -1, idxIfItHas = -1, idxNoteFirst = -1, idxNoteSecond = -1,\n    idxMSMap = -1, idxADT = -1, idxConvert = -1, idxLastEmpty = -1;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Run a test to make SURE\") !== -1) idxRunTest = i;\n  else if (t.indexOf(\"If it has, just call\") !== -1) idxIfItHas = i;\n  else if (t.indexOf(\"NOTE\") !== -1 && t.indexOf(\"has capacity for error messages\") !== -1) {\n    if (idxNoteFirst === -1) idxNoteFirst = i; else idxNoteSecond = i;\n  } else if (t === \"MSMap\") idxMSMap = i;\n  else if (t.indexOf(\"ADT that will create a map\") !== -1) idxADT = i;\n  else if (t.indexOf(\"convertToChar\") !== -1) idxConvert = i;\n}\n// The final (empty) paragraph of the document body.\nidxLastEmpty = paragraphs.items.length - 1;\n\n// 1) \"Run a test to make SURE ... check).\" -- collapse the run that was\n//    split around a <w:proofErr> \"double\" grammar-check marker back into a\n//    single run with identical text.\n{\n  const r = paragraphs.items[idxRunTest].getRange();\n  const ooxml = wrapBody(\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>Run a test to make SURE that the spot that was selected is legal. (double safe \\u2013 just check that it wasn\\u2019t already selected before \\u2013 shouldn\\u2019t be possible but never hurts to check).</w:t></w:r>' +\n    '</w:p>');\n  r.insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\n// 2) \"If it has, just call gui ... message?\" -- same kind of collapse\n//    around the \"gui\" spell/grammar-check markers.\n{\n  const r = paragraphs.items[idxIfItHas].getRange();\n  const ooxml = wrapBody(\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>If it has, just call gui \\u2013 maybe include an error message?</w:t></w:r>' +\n    '</w:p>');\n  r.insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\n// 3) First \"NOTE \\u2013 make it so that the gui has capacity ...\" -- collapse.\n{\n  const r = paragraphs.items[idxNoteFirst].getRange();\n  const ooxml = wrapBody(\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">NOTE \\u2013 make it so that the gui has capacity for error messages \\u2013 an option string passed into it. </w:t></w:r>' +\n    '</w:p>');\n  r.insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\n// 4) Remove the _GoBack bookmark from its old spot inside the \"ADT that\n//    will create a map...\" paragraph (it gets relocated in step 5).\n{\n  const r = paragraphs.items[idxADT].getRange();\n  const ooxml = wrapBody(\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>ADT that will create a map for minesweeper \\u2013 meaning its bomb and number locations. It stores them in itself as object type CellHold \\u2013 which will have either Bomb, Number, Blank set to true</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> \\u2013 or have 9 </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">be a bomb, 0 be a blank and 1-8 be the number </w:t></w:r>' +\n    '</w:p>');\n  r.insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\n// 5) \"MSMap\" heading -- relocate the _GoBack bookmark to sit inside this\n//    word, splitting it into \"MSM\" + bookmark + \"ap\".\n{\n  const r = paragraphs.items[idxMSMap].getRange();\n  const ooxml = wrapBody(\n    '<w:p><w:r><w:lastRenderedPageBreak/><w:t>MSM</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r><w:t>ap</w:t></w:r></w:p>');\n  r.insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\n// 6) \"Includes a function called convertToChar ...\" -- reflow the runs\n//    that were split around <w:proofErr> gramStart/gramEnd markers.\n{\n  const r = paragraphs.items[idxConvert].getRange();\n  const ooxml = wrapBody(\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Includes a function called convertToChar which will convert these numbers into their corresponding character values (\\u201cb\\u201d or \\u201c </w:t></w:r>' +\n    '<w:r><w:t>\\u201c or</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> \\u201c1-9\\u201d).</w:t></w:r>' +\n    '</w:p>');\n  r.insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\n// 7) Second \"NOTE \\u2013 make it so that the gui has capacity ...\" -- collapse.\n{\n  const r = paragraphs.items[idxNoteSecond].getRange();\n  const ooxml = wrapBody(\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">NOTE \\u2013 make it so that the gui has capacity for error messages \\u2013 an option string passed into it. </w:t></w:r>' +\n    '</w:p>');\n  r.insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\n// 8) The final, previously-empty bullet paragraph at the end of the\n//    document gains the new closing thought about an abstract class.\n{\n  const lastPara = paragraphs.items[idxLastEmpty];\n  lastPara.insertText(\n    \"Might want to use an abstract class to act as an interface, this will let me have multiple test classes that I can just extend into other testing files\\u2026. Not sure WHEN il make this but I def should. Maybe il make it after im done the first gui, then just pull stuff out as an abstract class. \",\n    Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "# Locate every paragraph whose text contains $needle; returns an array of\n# 1-based Paragraphs() indices (possibly empty, possibly more than one).\nfunction Find-AllParaIndices {\n    param($doc, $needle)\n    $n = $doc.Paragraphs.Count\n    $result = @()\n    for ($i = 1; $i -le $n; $i++) {\n        $t = $doc.Paragraphs($i).Range.Text\n        if ($t.IndexOf($needle) -ge 0) {\n            $result += $i\n        }\n    }\n    return $result\n}\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) \"Run a test to make SURE ... check).\" -- the run that used to be\n#    split around a <w:proofErr> \"double\" grammar-check marker collapses\n#    back into a single run with identical text.\n# ---------------------------------------------------------------------\n$idx = (Find-AllParaIndices $d \"Run a test to make SURE\")[0]\n$p = $d.Paragraphs($idx)\n$r = $p.Range\n$r.End = $r.End - 1\n$r.Text = \"Run a test to make SURE that the spot that was selected is legal. (double safe \u2013 just check that it wasn\u2019t already selected before \u2013 shouldn\u2019t be possible but never hurts to check).\"\n\n# ---------------------------------------------------------------------\n# 2) \"If it has, just call gui ... message?\" -- same kind of collapse\n#    around the \"gui\" spell/grammar-check markers.\n# ---------------------------------------------------------------------\n$idx = (Find-AllParaIndices $d \"If it has, just call\")[0]\n$p = $d.Paragraphs($idx)\n$r = $p.Range\n$r.End = $r.End - 1\n$r.Text = \"If it has, just call gui \u2013 maybe include an error message?\"\n\n# ---------------------------------------------------------------------\n# 3) Both \"NOTE \u2013 make it so that the gui has capacity ...\" paragraphs\n#    (there are two, identical, occurrences) -- collapse each.\n# ---------------------------------------------------------------------\n$noteIdxs = Find-AllParaIndices $d \"NOTE\"\nforeach ($idx in $noteIdxs) {\n    $p = $d.Paragraphs($idx)\n    $r = $p.Range\n    $r.End = $r.End - 1\n    $r.Text = \"NOTE \u2013 make it so that the gui has capacity for error messages \u2013 an option string passed into it. \"\n}\n\n# ---------------------------------------------------------------------\n# 4) Remove the _GoBack bookmark from its old spot inside the \"ADT that\n#    will create a map...\" paragraph (it gets relocated in step 5).\n# ---------------------------------------------------------------------\n$idx = (Find-AllParaIndices $d \"ADT that will create a map\")[0]\n$p = $d.Paragraphs($idx)\n$r = $p.Range\n$r.End = $r.End - 1\n$r.Delete()\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t>ADT that will create a map for minesweeper \u2013 meaning its bomb and number locations. It stores them in itself as object type CellHold \u2013 which will have either Bomb, Number, Blank set to true</w:t></w:r><w:r><w:t xml:space=\"preserve\"> \u2013 or have 9 </w:t></w:r><w:r><w:t xml:space=\"preserve\">be a bomb, 0 be a blank and 1-8 be the number </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$r.InsertXML($xml)\n\n# ---------------------------------------------------------------------\n# 5) \"MSMap\" heading -- relocate the _GoBack bookmark to sit inside this\n#    word, splitting it into \"MSM\" + bookmark + \"ap\".\n# ---------------------------------------------------------------------\n$idx = -1\n$n = $d.Paragraphs.Count\nfor ($i = 1; $i -le $n; $i++) {\n    if ($d.Paragraphs($i).Range.Text -eq \"MSMap`r\") { $idx = $i }\n}\n$p = $d.Paragraphs($idx)\n$r = $p.Range\n$r.End = $r.End - 1\n$r.Delete()\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>MSM</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t>ap</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$r.InsertXML($xml)\n\n# ---------------------------------------------------------------------\n# 6) \"Includes a function called convertToChar ...\" -- reflow the runs\n#    that were split around <w:proofErr> gramStart/gramEnd markers.\n# ---------------------------------------------------------------------\n$idx = (Find-AllParaIndices $d \"convertToChar\")[0]\n$p = $d.Paragraphs($idx)\n$r = $p.Range\n$r.End = $r.End - 1\n$r.Delete()\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Includes a function called convertToChar which will convert these numbers into their corresponding character values (\u201cb\u201d or \u201c </w:t></w:r><w:r><w:t>\u201c or</w:t></w:r><w:r><w:t xml:space=\"preserve\"> \u201c1-9\u201d).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$r.InsertXML($xml)\n\n# ---------------------------------------------------------------------\n# 7) The final, previously-empty bullet paragraph at the end of the\n#    document gains the new closing thought about an abstract class.\n# ---------------------------------------------------------------------\n$n = $d.Paragraphs.Count\n$last = $d.Paragraphs($n)\n$r = $last.Range\n$r.End = $r.End - 1\n$r.InsertBefore(\"Might want to use an abstract class to act as an interface, this will let me have multiple test classes that I can just extend into other testing files\u2026. Not sure WHEN il make this but I def should. Maybe il make it after im done the first gui, then just pull stuff out as an abstract class. \")\n"}
